$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting (bold font + box border + centered/top alignment) used by
# B1:H1 -- copy that format from H1 so the same style entry is reused.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I0 / IF columns, rows 2-69
$iVals = @(9,10,7,7,8,8,7,8,8,8,9,12,8,9,6,8,7,7,8,7,6,8,7,9,9,9,9,9,8,7,9,9,9,10,7,8,8,8,9,7,7,8,8,8,8,7,7,7,9,8,8,8,7,8,9,8,8,9,7,7,8,7,7,9,7,7,5,6)
$jVals = @(9,10,7,7,8,8,7,8,8,8,9,12,8,9,6,8,7,7,8,7,6,8,7,9,9,9,9,9,8,7,9,9,9,10,7,8,8,8,9,7,7,8,8,8,8,7,7,7,9,8,8,8,7,8,9,8,8,9,7,7,8,8,7,9,7,7,5,6)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
